# Refresh currentAveragePrice / LevePrice / LeveProfit columns (H:N) from the
# latest market-board pull across every Gathering/Crafting job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 158795.84
$ws.Range("J17").Value = 161290.02
$ws.Range("L17").Value = 483870.0599999999
$ws.Range("N17").Value = -484206.0599999999

# Row 33: Glazed and Confused
$ws.Range("H33").Value = 486.1579
$ws.Range("J33").Value = 825.5
$ws.Range("L33").Value = 825.5
$ws.Range("N33").Value = -1283.5

# Row 49: Going Nowhere Fast
$ws.Range("H49").Value = 405.66666
$ws.Range("I49").Value = 405.66666
$ws.Range("K49").Value = 1216.99998
$ws.Range("M49").Value = -1080.99998

# Row 126: Rebuilding to Code
$ws.Range("H126").Value = 94999.336
$ws.Range("J126").Value = 94999.336
$ws.Range("L126").Value = 94999.336
$ws.Range("N126").Value = -104879.336

# Row 135: For Tired Minds
$ws.Range("H135").Value = 1658.8462
$ws.Range("J135").Value = 7888
$ws.Range("L135").Value = 70992
$ws.Range("N135").Value = -76062

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1919.5667
$ws.Range("I137").Value = 2124.9412
$ws.Range("J137").Value = 1651
$ws.Range("K137").Value = 6374.823600000001
$ws.Range("L137").Value = 4953
$ws.Range("M137").Value = -3824.823600000001
$ws.Range("N137").Value = -10053

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 12746.148
$ws.Range("I32").Value = 13044.077
$ws.Range("K32").Value = 13044.077
$ws.Range("M32").Value = -12757.077

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 2531.0527
$ws.Range("I45").Value = 2057.5386
$ws.Range("K45").Value = 2057.5386
$ws.Range("M45").Value = -1680.5386

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 4746.7427
$ws.Range("I61").Value = 3054.204
$ws.Range("K61").Value = 3054.204
$ws.Range("M61").Value = -2842.204

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 3364.8708
$ws.Range("I74").Value = 2693.077
$ws.Range("K74").Value = 2693.077
$ws.Range("M74").Value = -1819.077

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 3364.8708
$ws.Range("I77").Value = 2693.077
$ws.Range("K77").Value = 13465.385
$ws.Range("M77").Value = -9097.385000000002

# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 22225630
$ws.Range("I102").Value = 2659.2
$ws.Range("K102").Value = 2659.2
$ws.Range("M102").Value = -1037.2

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 2540.5227
$ws.Range("I122").Value = 2104.0293
$ws.Range("J122").Value = 4024.6
$ws.Range("K122").Value = 6312.0879
$ws.Range("L122").Value = 12073.8
$ws.Range("M122").Value = -3862.0879
$ws.Range("N122").Value = -16973.8

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2219.103
$ws.Range("I132").Value = 2178.0508
$ws.Range("J132").Value = 2488.2222
$ws.Range("K132").Value = 6534.1524
$ws.Range("L132").Value = 7464.6666
$ws.Range("M132").Value = -4004.1524
$ws.Range("N132").Value = -12524.6666

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 4746.7427
$ws.Range("I136").Value = 3054.204
$ws.Range("K136").Value = 9162.612000000001
$ws.Range("M136").Value = -6612.612000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 2842.7144
$ws.Range("I20").Value = 2592.3845
$ws.Range("J20").Value = 3249.5
$ws.Range("K20").Value = 2592.3845
$ws.Range("L20").Value = 3249.5
$ws.Range("M20").Value = -2345.3845
$ws.Range("N20").Value = -3743.5

# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 4169291
$ws.Range("I86").Value = 6669716
$ws.Range("J86").Value = 1915.8334
$ws.Range("K86").Value = 6669716
$ws.Range("L86").Value = 1915.8334
$ws.Range("M86").Value = -6668593
$ws.Range("N86").Value = -4161.8334

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 4169291
$ws.Range("I89").Value = 6669716
$ws.Range("J89").Value = 1915.8334
$ws.Range("K89").Value = 33348580
$ws.Range("L89").Value = 9579.166999999999
$ws.Range("M89").Value = -33342964
$ws.Range("N89").Value = -20811.167

# Row 107: The Gold Experience
$ws.Range("H107").Value = 981.875
$ws.Range("I107").Value = 734.5714
$ws.Range("K107").Value = 734.5714
$ws.Range("M107").Value = 1185.4286

# Row 130: Annals of the Empire I
$ws.Range("H130").Value = 84974.75
$ws.Range("J130").Value = 84974.75
$ws.Range("L130").Value = 84974.75
$ws.Range("N130").Value = -95014.75

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 7598.659
$ws.Range("I134").Value = 2810.0715
$ws.Range("K134").Value = 8430.2145
$ws.Range("M134").Value = -5895.2145

$ws = $wb.Worksheets.Item("CRP")
# Row 23: Nothing to Hide
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

# Row 27: Behind the Mask
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

# Row 86: Birch, Please
$ws.Range("H86").Value = 5284.4443
$ws.Range("I86").Value = 3900.2
$ws.Range("J86").Value = 7014.75
$ws.Range("K86").Value = 3900.2
$ws.Range("L86").Value = 7014.75
$ws.Range("M86").Value = -2777.2
$ws.Range("N86").Value = -9260.75

# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 5284.4443
$ws.Range("I89").Value = 3900.2
$ws.Range("J89").Value = 7014.75
$ws.Range("K89").Value = 19501
$ws.Range("L89").Value = 35073.75
$ws.Range("M89").Value = -13885
$ws.Range("N89").Value = -46305.75

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 3091.8103
$ws.Range("I134").Value = 1333.175
$ws.Range("K134").Value = 3999.525
$ws.Range("M134").Value = -1464.525

$ws = $wb.Worksheets.Item("CUL")
# Row 62: Little Orphan Candy
$ws.Range("H62").Value = 8405
$ws.Range("I62").Value = 7756.5
$ws.Range("J62").Value = 10999
$ws.Range("K62").Value = 23269.5
$ws.Range("L62").Value = 32997
$ws.Range("M62").Value = -22583.5
$ws.Range("N62").Value = -34369

# Row 65: Confections of Confession (L)
$ws.Range("H65").Value = 8405
$ws.Range("I65").Value = 7756.5
$ws.Range("J65").Value = 10999
$ws.Range("K65").Value = 69808.5
$ws.Range("L65").Value = 98991
$ws.Range("M65").Value = -66376.5
$ws.Range("N65").Value = -105855

# Row 121: A Cookie for Your Troubles
$ws.Range("H121").Value = 3700.3
$ws.Range("I121").Value = 350
$ws.Range("J121").Value = 5933.8335
$ws.Range("K121").Value = 1050
$ws.Range("L121").Value = 17801.5005
$ws.Range("M121").Value = 260
$ws.Range("N121").Value = -20421.5005

# Row 134: Don't Knock It Till You've Tried It
$ws.Range("H134").Value = 2676.9167
$ws.Range("I134").Value = 1812.3
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 5436.9
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -366.8999999999996
$ws.Range("N134").Value = -31140

# Row 137: Creative Chocolate
$ws.Range("H137").Value = 3661.6428
$ws.Range("I137").Value = 3976.6667
$ws.Range("J137").Value = 3575.7273
$ws.Range("K137").Value = 11930.0001
$ws.Range("L137").Value = 10727.1819
$ws.Range("M137").Value = -6830.000100000001
$ws.Range("N137").Value = -20927.1819

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 8400
$ws.Range("I70").Value = 7800
$ws.Range("K70").Value = 7800
$ws.Range("M70").Value = -7530

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 8400
$ws.Range("I73").Value = 7800
$ws.Range("K73").Value = 7800
$ws.Range("M73").Value = -6864

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 3131.6667
$ws.Range("I126").Value = 2700
$ws.Range("J126").Value = 3995
$ws.Range("K126").Value = 8100
$ws.Range("L126").Value = 11985
$ws.Range("M126").Value = -5630
$ws.Range("N126").Value = -16925

$ws = $wb.Worksheets.Item("LTW")
# Row 29: Hands On
$ws.Range("H29").Value = 21333
$ws.Range("J29").Value = 21333
$ws.Range("L29").Value = 21333
$ws.Range("N29").Value = -21923

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 36192.75
$ws.Range("I61").Value = 38272.266
$ws.Range("K61").Value = 38272.266
$ws.Range("M61").Value = -38070.266

# Row 113: Peace in Rest
$ws.Range("H113").Value = 36192.75
$ws.Range("I113").Value = 38272.266
$ws.Range("K113").Value = 38272.266
$ws.Range("M113").Value = -36102.266

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 3161.3257
$ws.Range("J132").Value = 3169.3333
$ws.Range("L132").Value = 9507.999899999999
$ws.Range("N132").Value = -14567.9999

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 4006.6
$ws.Range("I136").Value = 3206.8462
$ws.Range("K136").Value = 9620.5386
$ws.Range("M136").Value = -7070.5386

$ws = $wb.Worksheets.Item("WVR")
# Row 100: Of Great Import
$ws.Range("H100").Value = 3832.8333
$ws.Range("J100").Value = 4998.6665
$ws.Range("L100").Value = 9997.333000000001
$ws.Range("N100").Value = -11079.333

# Row 113: A Tender Table
$ws.Range("H113").Value = 1084.6666
$ws.Range("I113").Value = 667.8
$ws.Range("J113").Value = 1293.1
$ws.Range("K113").Value = 2003.4
$ws.Range("L113").Value = 3879.3
$ws.Range("M113").Value = 166.6000000000001
$ws.Range("N113").Value = -8219.299999999999

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 2739.9048
$ws.Range("I132").Value = 2944.2
$ws.Range("K132").Value = 8832.599999999999
$ws.Range("M132").Value = -6302.599999999999

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 7913.7856
$ws.Range("J136").Value = 1433.3334
$ws.Range("L136").Value = 4300.0002
$ws.Range("N136").Value = -9400.0002

Write-Output "Sheets refreshed."
